# Auto update: 2025-11-29 04:59:05
#
# The daily refresh swapped which row holds "HD HYUNDAI MIPO" (010620.KS)
# vs "HDKSOE" (009540.KS) and refreshed their price/score metrics, plus
# updated the shared MACRO_SCORE (column N) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 becomes HD HYUNDAI MIPO (010620.KS) with refreshed metrics
$ws.Range("B4").Value = "HD HYUNDAI MIPO"
$ws.Range("C4").Value = "010620.KS"
$ws.Range("D4").Value = 223000
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 20
$ws.Range("H4").Value = 46
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 63
$ws.Range("K4").Value = 51.8

# Row 5 becomes HDKSOE (009540.KS) with refreshed metrics
$ws.Range("B5").Value = "HDKSOE"
$ws.Range("C5").Value = "009540.KS"
$ws.Range("D5").Value = 410000
$ws.Range("F5").Value = -2.38
$ws.Range("G5").Value = 10
$ws.Range("H5").Value = 70
$ws.Range("I5").Value = 56
$ws.Range("J5").Value = 80
$ws.Range("K5").Value = 51.2

# MACRO_SCORE (column N) refreshed for every data row
$ws.Range("N2").Value = 85.8724807945396
$ws.Range("N3").Value = 85.8724807945396
$ws.Range("N4").Value = 85.8724807945396
$ws.Range("N5").Value = 85.8724807945396
